$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 423.2
$ws.Range("I2").Value = 294.57144
$ws.Range("J2").Value = 723.3333
$ws.Range("K2").Value = 294.57144
$ws.Range("L2").Value = 723.3333
$ws.Range("M2").Value = -181.57144
$ws.Range("N2").Value = -949.3333
$ws.Range("H26").Value = 9999
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 9999
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 9999
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -10687
$ws.Range("H51").Value = 4867.2915
$ws.Range("J51").Value = 4934.3477
$ws.Range("L51").Value = 4934.3477
$ws.Range("N51").Value = -5902.3477
$ws.Range("H113").Value = 4880.5835
$ws.Range("J113").Value = 5406.8
$ws.Range("L113").Value = 5406.8
$ws.Range("N113").Value = -11914.8
$ws.Range("H131").Value = 5564.8667
$ws.Range("I131").Value = 3421
$ws.Range("K131").Value = 10263
$ws.Range("M131").Value = -5223
$ws.Range("H135").Value = 783.9355
$ws.Range("I135").Value = 539.0909
$ws.Range("K135").Value = 4851.8181
$ws.Range("M135").Value = -2316.8181

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 35
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 35
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 35
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -267
$ws.Range("H8").Value = 3000000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 3000000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 3000000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -3000288
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H11").Value = 10000000
$ws.Range("I11").Value = 10000000
$ws.Range("K11").Value = 10000000
$ws.Range("M11").Value = -9999856
$ws.Range("H45").Value = 1908
$ws.Range("I45").Value = 1640
$ws.Range("J45").Value = 2354.6667
$ws.Range("K45").Value = 1640
$ws.Range("L45").Value = 2354.6667
$ws.Range("M45").Value = -1263
$ws.Range("N45").Value = -3108.6667
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H102").Value = 145028.58
$ws.Range("I102").Value = 201679.7
$ws.Range("K102").Value = 201679.7
$ws.Range("M102").Value = -200057.7
$ws.Range("H122").Value = 2780.394
$ws.Range("I122").Value = 1832.7222
$ws.Range("K122").Value = 5498.1666
$ws.Range("M122").Value = -3048.1666

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14293.058
$ws.Range("I86").Value = 6929.0835
$ws.Range("J86").Value = 30862
$ws.Range("K86").Value = 6929.0835
$ws.Range("L86").Value = 30862
$ws.Range("M86").Value = -5806.0835
$ws.Range("N86").Value = -33108
$ws.Range("H89").Value = 14293.058
$ws.Range("I89").Value = 6929.0835
$ws.Range("J89").Value = 30862
$ws.Range("K89").Value = 34645.4175
$ws.Range("L89").Value = 154310
$ws.Range("M89").Value = -29029.4175
$ws.Range("N89").Value = -165542
$ws.Range("H105").Value = 21239.8
$ws.Range("I105").Value = 34283
$ws.Range("K105").Value = 34283
$ws.Range("M105").Value = -32536
$ws.Range("H107").Value = 3809.5334
$ws.Range("J107").Value = 3374.75
$ws.Range("L107").Value = 3374.75
$ws.Range("N107").Value = -7214.75

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 8084.6924
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1200
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H56").Value = 12995
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 12995
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 12995
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -14685

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1439
$ws.Range("J23").Value = 1758.6666
$ws.Range("L23").Value = 5275.9998
$ws.Range("N23").Value = -5745.9998
$ws.Range("H133").Value = 4654.5
$ws.Range("I133").Value = 4654.5
$ws.Range("K133").Value = 13963.5
$ws.Range("M133").Value = -8903.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 403980.6
$ws.Range("I70").Value = 503726
$ws.Range("J70").Value = 4999
$ws.Range("K70").Value = 503726
$ws.Range("L70").Value = 4999
$ws.Range("M70").Value = -503456
$ws.Range("N70").Value = -5539
$ws.Range("H73").Value = 403980.6
$ws.Range("I73").Value = 503726
$ws.Range("J73").Value = 4999
$ws.Range("K73").Value = 503726
$ws.Range("L73").Value = 4999
$ws.Range("M73").Value = -502790
$ws.Range("N73").Value = -6871
$ws.Range("H80").Value = 3321.5
$ws.Range("I80").Value = 3474.5
$ws.Range("K80").Value = 3474.5
$ws.Range("M80").Value = -2476.5
$ws.Range("H83").Value = 3321.5
$ws.Range("I83").Value = 3474.5
$ws.Range("K83").Value = 17372.5
$ws.Range("M83").Value = -12380.5
$ws.Range("H128").Value = 65000
$ws.Range("J128").Value = 65000
$ws.Range("L128").Value = 65000
$ws.Range("N128").Value = -74960

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 212.14285
$ws.Range("I9").Value = 239.16667
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 239.16667
$ws.Range("L9").Value = 50
$ws.Range("M9").Value = -15.16667000000001
$ws.Range("N9").Value = -498
$ws.Range("H22").Value = 1895.1818
$ws.Range("I22").Value = 524.6667
$ws.Range("J22").Value = 3539.8
$ws.Range("K22").Value = 524.6667
$ws.Range("L22").Value = 3539.8
$ws.Range("M22").Value = -229.6667
$ws.Range("N22").Value = -4129.8
$ws.Range("H27").Value = 1895.1818
$ws.Range("I27").Value = 524.6667
$ws.Range("J27").Value = 3539.8
$ws.Range("K27").Value = 524.6667
$ws.Range("L27").Value = 3539.8
$ws.Range("M27").Value = -417.6667
$ws.Range("N27").Value = -3753.8
$ws.Range("H55").Value = 651.13043
$ws.Range("I55").Value = 430.5
$ws.Range("J55").Value = 994.3333
$ws.Range("K55").Value = 430.5
$ws.Range("L55").Value = 994.3333
$ws.Range("M55").Value = -257.5
$ws.Range("N55").Value = -1340.3333
$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H74").Value = 57497.25
$ws.Range("I74").Value = 55995
$ws.Range("K74").Value = 55995
$ws.Range("M74").Value = -54997
$ws.Range("H77").Value = 57497.25
$ws.Range("I77").Value = 55995
$ws.Range("K77").Value = 167985
$ws.Range("M77").Value = -162993
$ws.Range("H87").Value = 52397.8
$ws.Range("J87").Value = 51997.5
$ws.Range("L87").Value = 51997.5
$ws.Range("N87").Value = -54243.5
$ws.Range("H90").Value = 52397.8
$ws.Range("J90").Value = 51997.5
$ws.Range("L90").Value = 155992.5
$ws.Range("N90").Value = -167224.5
$ws.Range("H122").Value = 4795.9697
$ws.Range("I122").Value = 3947.1667
$ws.Range("J122").Value = 5281
$ws.Range("K122").Value = 11841.5001
$ws.Range("L122").Value = 15843
$ws.Range("M122").Value = -9391.500100000001
$ws.Range("N122").Value = -20743

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 3570.2856
$ws.Range("I13").Value = 1748.25
$ws.Range("J13").Value = 5999.6665
$ws.Range("K13").Value = 1748.25
$ws.Range("L13").Value = 5999.6665
$ws.Range("M13").Value = -1608.25
$ws.Range("N13").Value = -6279.6665
$ws.Range("H57").Value = 67492.5
$ws.Range("J57").Value = 69990
$ws.Range("L57").Value = 69990
$ws.Range("N57").Value = -71498
$ws.Range("H74").Value = 44323.332
$ws.Range("J74").Value = 44323.332
$ws.Range("L74").Value = 44323.332
$ws.Range("N74").Value = -46195.332
$ws.Range("H77").Value = 44323.332
$ws.Range("J77").Value = 44323.332
$ws.Range("L77").Value = 132969.996
$ws.Range("N77").Value = -142329.996
